$wb = $excel.ActiveWorkbook

# Update values on "yeni_degiskenler" sheet
$ws = $wb.Worksheets.Item("yeni_degiskenler")
$ws.Range("B3").Value = 0.25
$ws.Range("B4").Value = 0.3

# Update the selected cell on that sheet to match the diff (B5 instead of K20)
$ws.Range("B5").Select()
